$wb = $excel.ActiveWorkbook

# --- Rename existing sheet and add the new "Raw Data" sheet right after it ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "FX History Report Result"
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Raw Data"

# --- Populate "FX History Report Result" with the FX table ---
$headers = @("Company","Currency","Code","Buy Rate","Sell Rate","Rate Last Modified","Active")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws1.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$ws1.Range("A2").Value = "Australia"
$ws1.Range("B2").Value = "AUD"
$ws1.Range("C2").Value = "22/05/2022"
$ws1.Range("D2").Value = 1
$ws1.Range("E2").Value = 1

$ws1.Range("A3").Value = "United Kingdom"
$ws1.Range("B3").Value = "GBP"
$ws1.Range("C3").Value = "21/05/2022"
$ws1.Range("D3").Value = 0.43
$ws1.Range("E3").Value = 0.47

$ws1.Range("A4").Value = "United States"
$ws1.Range("B4").Value = "USD"
$ws1.Range("C4").Value = "20/05/2022"
$ws1.Range("D4").Value = 0.75
$ws1.Range("E4").Value = 0.75

$ws1.Range("A5").Value = "South African"
$ws1.Range("B5").Value = "ZAR"
$ws1.Range("C5").Value = "19/05/2022"
$ws1.Range("D5").Value = 2.35
$ws1.Range("E5").Value = 2.38

# --- Left-align only the populated cells (creates the new cellXfs entry
#     without materialising style-only entries for cells that stay empty) ---
$ws1.Range("A1:G1").HorizontalAlignment = -4131
$ws1.Range("A2:E2").HorizontalAlignment = -4131
$ws1.Range("A3:E3").HorizontalAlignment = -4131
$ws1.Range("A4:E4").HorizontalAlignment = -4131
$ws1.Range("A5:E5").HorizontalAlignment = -4131

# --- Best-fit the columns like the original report ---
$ws1.Columns.Item(1).AutoFit() | Out-Null
$ws1.Columns.Item(2).AutoFit() | Out-Null
$ws1.Columns.Item(3).AutoFit() | Out-Null
$ws1.Columns.Item(4).AutoFit() | Out-Null
$ws1.Columns.Item(5).AutoFit() | Out-Null
$ws1.Columns.Item(6).AutoFit() | Out-Null
$ws1.Columns.Item(7).AutoFit() | Out-Null

# --- Page setup (portrait) to match the exported report ---
$ws1.PageSetup.Orientation = 1

# --- Selection state observed in the authored file ---
[void]$ws1.Range("C8").Select()
$ws1.Select()
